$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 47.837127
$ws.Range("H2").Value = 143.511381
$ws.Range("I2").Value = 0.05107837591710958
$ws.Range("J2").Value = 0.05107837591710957
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 163.7119853333333
$ws.Range("N2").Value = 491.135956
$ws.Range("O2").Value = 0.2754003062401033
$ws.Range("P2").Value = 0.2754003062401033
$ws.Range("Q2").Value = 7831.511033812804
$ws.Range("R2").Value = 70483.59930431523
$ws.Range("S2").Value = 0.0140670003698191
$ws.Range("T2").Value = 0.01406700036981909

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 47.837127
$ws.Range("H3").Value = 143.511381
$ws.Range("I3").Value = 0.05107837591710958
$ws.Range("J3").Value = 0.05107837591710957
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 172.558497
$ws.Range("N3").Value = 517.675491
$ws.Range("O3").Value = 0.290282124557779
$ws.Range("P3").Value = 0.290282124557779
$ws.Range("Q3").Value = 8254.70273591812
$ws.Range("R3").Value = 74292.32462326306
$ws.Range("S3").Value = 0.01482713948017946
$ws.Range("T3").Value = 0.01482713948017946

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 47.837127
$ws.Range("H4").Value = 143.511381
$ws.Range("I4").Value = 0.05107837591710958
$ws.Range("J4").Value = 0.05107837591710957
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 111.4881643333333
$ws.Range("N4").Value = 334.464493
$ws.Range("O4").Value = 0.1875481171218523
$ws.Range("P4").Value = 0.1875481171218523
$ws.Range("Q4").Value = 5333.273476210537
$ws.Range("R4").Value = 47999.46128589483
$ws.Range("S4").Value = 0.009579653228896068
$ws.Range("T4").Value = 0.009579653228896068

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 47.837127
$ws.Range("H5").Value = 143.511381
$ws.Range("I5").Value = 0.05107837591710958
$ws.Range("J5").Value = 0.05107837591710957
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 86.95798233333333
$ws.Range("N5").Value = 260.873947
$ws.Range("O5").Value = 0.1462828449356383
$ws.Range("P5").Value = 0.1462828449356383
$ws.Range("Q5").Value = 4159.820044543423
$ws.Range("R5").Value = 37438.3804008908
$ws.Range("S5").Value = 0.007471890143846785
$ws.Range("T5").Value = 0.007471890143846784

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 47.837127
$ws.Range("H6").Value = 143.511381
$ws.Range("I6").Value = 0.05107837591710958
$ws.Range("J6").Value = 0.05107837591710957
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 59.73436333333333
$ws.Range("N6").Value = 179.20309
$ws.Range("O6").Value = 0.100486607144627
$ws.Range("P6").Value = 0.100486607144627
$ws.Range("Q6").Value = 2857.52032504081
$ws.Range("R6").Value = 25717.68292536729
$ws.Range("S6").Value = 0.005132692694368166
$ws.Range("T6").Value = 0.005132692694368165

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 255.905248
$ws.Range("H7").Value = 767.715744
$ws.Range("I7").Value = 0.2732443454747012
$ws.Range("J7").Value = 0.2732443454747011
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 163.7119853333333
$ws.Range("N7").Value = 491.135956
$ws.Range("O7").Value = 0.2754003062401033
$ws.Range("P7").Value = 0.2754003062401033
$ws.Range("Q7").Value = 41894.75620729903
$ws.Range("R7").Value = 377052.8058656912
$ws.Range("S7").Value = 0.07525157642210929
$ws.Range("T7").Value = 0.07525157642210928

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 255.905248
$ws.Range("H8").Value = 767.715744
$ws.Range("I8").Value = 0.2732443454747012
$ws.Range("J8").Value = 0.2732443454747011
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 172.558497
$ws.Range("N8").Value = 517.675491
$ws.Range("O8").Value = 0.290282124557779
$ws.Range("P8").Value = 0.290282124557779
$ws.Range("Q8").Value = 44158.62496929225
$ws.Range("R8").Value = 397427.6247236303
$ws.Range("S8").Value = 0.079317949127796
$ws.Range("T8").Value = 0.07931794912779598

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 255.905248
$ws.Range("H9").Value = 767.715744
$ws.Range("I9").Value = 0.2732443454747012
$ws.Range("J9").Value = 0.2732443454747011
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 111.4881643333333
$ws.Range("N9").Value = 334.464493
$ws.Range("O9").Value = 0.1875481171218523
$ws.Range("P9").Value = 0.1875481171218523
$ws.Range("Q9").Value = 28530.40634278642
$ws.Range("R9").Value = 256773.6570850778
$ws.Range("S9").Value = 0.05124646250797312
$ws.Range("T9").Value = 0.05124646250797312

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 255.905248
$ws.Range("H10").Value = 767.715744
$ws.Range("I10").Value = 0.2732443454747012
$ws.Range("J10").Value = 0.2732443454747011
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 86.95798233333333
$ws.Range("N10").Value = 260.873947
$ws.Range("O10").Value = 0.1462828449356383
$ws.Range("P10").Value = 0.1462828449356383
$ws.Range("Q10").Value = 22253.00403459129
$ws.Range("R10").Value = 200277.0363113215
$ws.Range("S10").Value = 0.0399709602186157
$ws.Range("T10").Value = 0.03997096021861569

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 255.905248
$ws.Range("H11").Value = 767.715744
$ws.Range("I11").Value = 0.2732443454747012
$ws.Range("J11").Value = 0.2732443454747011
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 59.73436333333333
$ws.Range("N11").Value = 179.20309
$ws.Range("O11").Value = 0.100486607144627
$ws.Range("P11").Value = 0.100486607144627
$ws.Range("Q11").Value = 15286.33706293877
$ws.Range("R11").Value = 137577.0335664489
$ws.Range("S11").Value = 0.02745739719820703
$ws.Range("T11").Value = 0.02745739719820702

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 295.7629596666666
$ws.Range("H12").Value = 887.288879
$ws.Range("I12").Value = 0.315802653370277
$ws.Range("J12").Value = 0.315802653370277
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 163.7119853333333
$ws.Range("N12").Value = 491.135956
$ws.Range("O12").Value = 0.2754003062401033
$ws.Range("P12").Value = 0.2754003062401033
$ws.Range("Q12").Value = 48419.94131509258
$ws.Range("R12").Value = 435779.4718358333
$ws.Range("S12").Value = 0.08697214744961147
$ws.Range("T12").Value = 0.08697214744961147

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 295.7629596666666
$ws.Range("H13").Value = 887.288879
$ws.Range("I13").Value = 0.315802653370277
$ws.Range("J13").Value = 0.315802653370277
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 172.558497
$ws.Range("N13").Value = 517.675491
$ws.Range("O13").Value = 0.290282124557779
$ws.Range("P13").Value = 0.290282124557779
$ws.Range("Q13").Value = 51036.41178835161
$ws.Range("R13").Value = 459327.7060951645
$ws.Range("S13").Value = 0.09167186516130783
$ws.Range("T13").Value = 0.09167186516130783

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 295.7629596666666
$ws.Range("H14").Value = 887.288879
$ws.Range("I14").Value = 0.315802653370277
$ws.Range("J14").Value = 0.315802653370277
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 111.4881643333333
$ws.Range("N14").Value = 334.464493
$ws.Range("O14").Value = 0.1875481171218523
$ws.Range("P14").Value = 0.1875481171218523
$ws.Range("Q14").Value = 32974.06945103037
$ws.Range("R14").Value = 296766.6250592733
$ws.Range("S14").Value = 0.05922819302168043
$ws.Range("T14").Value = 0.05922819302168043

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 295.7629596666666
$ws.Range("H15").Value = 887.288879
$ws.Range("I15").Value = 0.315802653370277
$ws.Range("J15").Value = 0.315802653370277
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 86.95798233333333
$ws.Range("N15").Value = 260.873947
$ws.Range("O15").Value = 0.1462828449356383
$ws.Range("P15").Value = 0.1462828449356383
$ws.Range("Q15").Value = 25718.95022154838
$ws.Range("R15").Value = 231470.5519939354
$ws.Range("S15").Value = 0.04619651057322737
$ws.Range("T15").Value = 0.04619651057322737

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 295.7629596666666
$ws.Range("H16").Value = 887.288879
$ws.Range("I16").Value = 0.315802653370277
$ws.Range("J16").Value = 0.315802653370277
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 59.73436333333333
$ws.Range("N16").Value = 179.20309
$ws.Range("O16").Value = 0.100486607144627
$ws.Range("P16").Value = 0.100486607144627
$ws.Range("Q16").Value = 17667.21209327067
$ws.Range("R16").Value = 159004.9088394361
$ws.Range("S16").Value = 0.03173393716444983
$ws.Range("T16").Value = 0.03173393716444983

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 150.2412563333333
$ws.Range("H17").Value = 450.723769
$ws.Range("I17").Value = 0.1604209920309976
$ws.Range("J17").Value = 0.1604209920309976
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 163.7119853333333
$ws.Range("N17").Value = 491.135956
$ws.Range("O17").Value = 0.2754003062401033
$ws.Range("P17").Value = 0.2754003062401033
$ws.Range("Q17").Value = 24596.29435330424
$ws.Range("R17").Value = 221366.6491797382
$ws.Range("S17").Value = 0.04417999033267792
$ws.Range("T17").Value = 0.04417999033267791

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 150.2412563333333
$ws.Range("H18").Value = 450.723769
$ws.Range("I18").Value = 0.1604209920309976
$ws.Range("J18").Value = 0.1604209920309976
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 172.558497
$ws.Range("N18").Value = 517.675491
$ws.Range("O18").Value = 0.290282124557779
$ws.Range("P18").Value = 0.290282124557779
$ws.Range("Q18").Value = 25925.40538027173
$ws.Range("R18").Value = 233328.6484224456
$ws.Range("S18").Value = 0.04656734639042451
$ws.Range("T18").Value = 0.04656734639042451

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 150.2412563333333
$ws.Range("H19").Value = 450.723769
$ws.Range("I19").Value = 0.1604209920309976
$ws.Range("J19").Value = 0.1604209920309976
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 111.4881643333333
$ws.Range("N19").Value = 334.464493
$ws.Range("O19").Value = 0.1875481171218523
$ws.Range("P19").Value = 0.1875481171218523
$ws.Range("Q19").Value = 16750.12187573712
$ws.Range("R19").Value = 150751.0968816341
$ws.Range("S19").Value = 0.03008665500223327
$ws.Range("T19").Value = 0.03008665500223327

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 150.2412563333333
$ws.Range("H20").Value = 450.723769
$ws.Range("I20").Value = 0.1604209920309976
$ws.Range("J20").Value = 0.1604209920309976
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 86.95798233333333
$ws.Range("N20").Value = 260.873947
$ws.Range("O20").Value = 0.1462828449356383
$ws.Range("P20").Value = 0.1462828449356383
$ws.Range("Q20").Value = 13064.6765139718
$ws.Range("R20").Value = 117582.0886257462
$ws.Range("S20").Value = 0.0234668391016917
$ws.Range("T20").Value = 0.02346683910169169

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 150.2412563333333
$ws.Range("H21").Value = 450.723769
$ws.Range("I21").Value = 0.1604209920309976
$ws.Range("J21").Value = 0.1604209920309976
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 59.73436333333333
$ws.Range("N21").Value = 179.20309
$ws.Range("O21").Value = 0.100486607144627
$ws.Range("P21").Value = 0.100486607144627
$ws.Range("Q21").Value = 8974.5657934718
$ws.Range("R21").Value = 80771.0921412462
$ws.Range("S21").Value = 0.01612016120397019
$ws.Range("T21").Value = 0.01612016120397019

$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 186.7970273333333
$ws.Range("H22").Value = 560.391082
$ws.Range("I22").Value = 0.1994536332069146
$ws.Range("J22").Value = 0.1994536332069146
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 163.7119853333333
$ws.Range("N22").Value = 491.135956
$ws.Range("O22").Value = 0.2754003062401033
$ws.Range("P22").Value = 0.2754003062401033
$ws.Range("Q22").Value = 30580.91219910493
$ws.Range("R22").Value = 275228.2097919444
$ws.Range("S22").Value = 0.05492959166588553
$ws.Range("T22").Value = 0.05492959166588553

$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 186.7970273333333
$ws.Range("H23").Value = 560.391082
$ws.Range("I23").Value = 0.1994536332069146
$ws.Range("J23").Value = 0.1994536332069146
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 172.558497
$ws.Range("N23").Value = 517.675491
$ws.Range("O23").Value = 0.290282124557779
$ws.Range("P23").Value = 0.290282124557779
$ws.Range("Q23").Value = 32233.41428070791
$ws.Range("R23").Value = 290100.7285263712
$ws.Range("S23").Value = 0.05789782439807115
$ws.Range("T23").Value = 0.05789782439807115

$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 186.7970273333333
$ws.Range("H24").Value = 560.391082
$ws.Range("I24").Value = 0.1994536332069146
$ws.Range("J24").Value = 0.1994536332069146
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 111.4881643333333
$ws.Range("N24").Value = 334.464493
$ws.Range("O24").Value = 0.1875481171218523
$ws.Range("P24").Value = 0.1875481171218523
$ws.Range("Q24").Value = 20825.65768031682
$ws.Range("R24").Value = 187430.9191228514
$ws.Range("S24").Value = 0.0374071533610694
$ws.Range("T24").Value = 0.0374071533610694

$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 186.7970273333333
$ws.Range("H25").Value = 560.391082
$ws.Range("I25").Value = 0.1994536332069146
$ws.Range("J25").Value = 0.1994536332069146
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 86.95798233333333
$ws.Range("N25").Value = 260.873947
$ws.Range("O25").Value = 0.1462828449356383
$ws.Range("P25").Value = 0.1462828449356383
$ws.Range("Q25").Value = 16243.49260277118
$ws.Range("R25").Value = 146191.4334249407
$ws.Range("S25").Value = 0.02917664489825678
$ws.Range("T25").Value = 0.02917664489825678

$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 186.7970273333333
$ws.Range("H26").Value = 560.391082
$ws.Range("I26").Value = 0.1994536332069146
$ws.Range("J26").Value = 0.1994536332069146
$ws.Range("K26").Value = 3
$ws.Range("M26").Value = 59.73436333333333
$ws.Range("N26").Value = 179.20309
$ws.Range("O26").Value = 0.100486607144627
$ws.Range("P26").Value = 0.100486607144627
$ws.Range("Q26").Value = 11158.20150031593
$ws.Range("R26").Value = 100423.8135028434
$ws.Range("S26").Value = 0.02004241888363176
$ws.Range("T26").Value = 0.02004241888363175

